$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# FAQS sheet: add two new greeting/farewell rows (12 and 13) with a bordered,
# wrap-text style (left/right medium #CCCCCC border) matching the style used
# for the header row's font.
# ---------------------------------------------------------------------------
$wsFaqs = $wb.Worksheets.Item("FAQS")

$wsFaqs.Range("A12").Value = "Hi"
$wsFaqs.Range("B12").Value = "Hello, how may I help you?"
$wsFaqs.Range("A13").Value = "Bye"
$wsFaqs.Range("B13").Value = "Goodbye, have a nice day."

# Build the new border/wrap style on A12 first ...
$a12 = $wsFaqs.Range("A12")
$a12.WrapText = $true
$a12.Borders.Item(7).Weight = -4138
$a12.Borders.Item(7).Color = 13421772
$a12.Borders.Item(10).Weight = -4138
$a12.Borders.Item(10).Color = 13421772

# ... then copy that exact format onto the other three new cells so they all
# share the same single new style entry instead of each minting their own.
$a12.Copy()
$wsFaqs.Range("B12").PasteSpecial(-4122)
$wsFaqs.Range("A13").PasteSpecial(-4122)
$wsFaqs.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# GENERAL_INTENTS sheet: add a new "deposit" intent row.
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("GENERAL_INTENTS")
$wsGeneral.Range("A5").Value = "I want to deposit AMOUNT"
$wsGeneral.Range("B5").Value = "AMOUNT has been deposited into"

# ---------------------------------------------------------------------------
# Selections / active sheet, matching where the author ended up afterwards.
# ---------------------------------------------------------------------------
$wsFaqs.Range("A14").Select()
$wsGeneral.Range("B5").Select()
